$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-number cells to use a text number format before assignment
# so Excel does not auto-convert numeric-looking strings into numbers.
$textCells = @('D5', 'D6', 'D10', 'D11', 'D12', 'D14', 'D18', 'D20', 'D21', 'D22', 'D23', 'D24', 'D27', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D36', 'D39', 'D40', 'D42', 'D43', 'D44', 'D48', 'D49', 'D50', 'D51')
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '64.077.40'
$ws.Range('E2').Value = '  -0.16%  '
$ws.Range('D3').Value = '2.758.90'
$ws.Range('E3').Value = '  +0.96%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '577.79'
$ws.Range('E5').Value = '  -0.62%  '
$ws.Range('D6').Value = '158.99'
$ws.Range('E6').Value = '  +0.72%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  -3.21%  '
$ws.Range('E9').Value = '  -1.81%  '
$ws.Range('D10').Value = '0.163'
$ws.Range('E10').Value = '  +2.18%  '
$ws.Range('D11').Value = '0.385'
$ws.Range('E11').Value = '  -2.57%  '
$ws.Range('D12').Value = '5.64'
$ws.Range('E12').Value = '  -16.77%  '
$ws.Range('D13').Value = '3.246.24'
$ws.Range('E13').Value = '  +0.25%  '
$ws.Range('D14').Value = '26.95'
$ws.Range('E14').Value = '  -1.46%  '
$ws.Range('D15').Value = '63.664.84'
$ws.Range('E15').Value = '  -0.65%  '
$ws.Range('E16').Value = '  -1.66%  '
$ws.Range('D17').Value = '2.763.41'
$ws.Range('E17').Value = '  +0.18%  '
$ws.Range('D18').Value = '12.21'
$ws.Range('E18').Value = '  +1.01%  '
$ws.Range('E19').Value = '  -1.83%  '
$ws.Range('D20').Value = '358.05'
$ws.Range('E20').Value = '  -1.55%  '
$ws.Range('D21').Value = '6.76'
$ws.Range('E21').Value = '  -3.17%  '
$ws.Range('B22').Value = 'Polygon'
$ws.Range('C22').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D22').Value = '0.536'
$ws.Range('E22').Value = '  -1.81%  '
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('D24').Value = '65.47'
$ws.Range('E24').Value = '  -2.51%  '
$ws.Range('E25').Value = '  -0.87%  '
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('D27').Value = '0.998'
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('E28').Value = '  -1.70%  '
$ws.Range('D29').Value = '1.97'
$ws.Range('E29').Value = '  -2.75%  '
$ws.Range('D30').Value = '7.29'
$ws.Range('E30').Value = '  +0.38%  '
$ws.Range('D31').Value = '1.27'
$ws.Range('E31').Value = '  +0.13%  '
$ws.Range('D32').Value = '169.93'
$ws.Range('E32').Value = '  -2.21%  '
$ws.Range('D33').Value = '4.98'
$ws.Range('E33').Value = '  +1.18%  '
$ws.Range('D34').Value = '20.27'
$ws.Range('E34').Value = '  -1.96%  '
$ws.Range('E35').Value = '  +1.07%  '
$ws.Range('D36').Value = '0.998'
$ws.Range('E36').Value = '  +0.07%  '
$ws.Range('E37').Value = '  -1.10%  '
$ws.Range('E38').Value = '  -0.65%  '
$ws.Range('D39').Value = '6.32'
$ws.Range('E39').Value = '  +2.51%  '
$ws.Range('D40').Value = '338.11'
$ws.Range('E40').Value = '  -0.17%  '
$ws.Range('E41').Value = '  -2.66%  '
$ws.Range('D42').Value = '39.17'
$ws.Range('E42').Value = '  -0.76%  '
$ws.Range('D43').Value = '21.59'
$ws.Range('E43').Value = '  -1.86%  '
$ws.Range('D44').Value = '21.86'
$ws.Range('E44').Value = '  -1.91%  '
$ws.Range('E45').Value = '  -2.49%  '
$ws.Range('E46').Value = '  -1.70%  '
$ws.Range('E47').Value = '  -0.26%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').Value = '0.631'
$ws.Range('E48').Value = '  -2.87%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').Value = '135.51'
$ws.Range('E49').Value = '  -1.61%  '
$ws.Range('D50').Value = '0.999'
$ws.Range('E50').Value = '  +0.00%  '
$ws.Range('D51').Value = '11.06'
$ws.Range('E51').Value = '  +0.12%  '

# Restore the default (Normal) style on cells where we had to force a text
# number format, so formatting matches the original workbook.
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = "Normal"
}

